# Apply cryptocurrency price/volume updates (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.707.73"
$ws.Range("E2").Value = "  +2.16%  "
$ws.Range("D3").Value = "3.372.82"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'592.83"
$ws.Range("E5").Value = "  +6.32%  "
$ws.Range("D6").Value = "'187.04"
$ws.Range("E6").Value = "  -1.42%  "
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "'0.600"
$ws.Range("E7").Value = "  +2.84%  "
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.184"
$ws.Range("E9").Value = "  +2.66%  "
$ws.Range("D10").Value = "'0.590"
$ws.Range("E10").Value = "  +1.40%  "
$ws.Range("D11").Value = "'47.44"
$ws.Range("E11").Value = "  +2.24%  "
$ws.Range("E12").Value = "  +3.26%  "
$ws.Range("D13").Value = "3.920.32"
$ws.Range("E13").Value = "  +1.09%  "
$ws.Range("D14").Value = "'637.55"
$ws.Range("E14").Value = "  +8.39%  "
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D16").Value = "67.791.54"
$ws.Range("E16").Value = "  +2.39%  "
$ws.Range("D17").Value = "3.378.35"
$ws.Range("E17").Value = "  +1.37%  "
$ws.Range("E18").Value = "  +1.26%  "
$ws.Range("D19").Value = "'18.08"
$ws.Range("E19").Value = "  +0.91%  "
$ws.Range("D20").Value = "'11.16"
$ws.Range("E20").Value = "  +1.28%  "
$ws.Range("D21").Value = "'0.910"
$ws.Range("E21").Value = "  +1.24%  "
$ws.Range("D22").Value = "'17.93"
$ws.Range("E22").Value = "  -1.26%  "
$ws.Range("D23").Value = "'5.11"
$ws.Range("E23").Value = "  +1.74%  "
$ws.Range("D24").Value = "'99.88"
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("D25").Value = "'4.04"
$ws.Range("E25").Value = "  +1.73%  "
$ws.Range("E26").Value = "  +5.56%  "
$ws.Range("D27").Value = "'9.77"
$ws.Range("E27").Value = "  +3.33%  "
$ws.Range("D28").Value = "'32.82"
$ws.Range("E28").Value = "  +6.69%  "
$ws.Range("D29").Value = "'8.72"
$ws.Range("E29").Value = "  +2.64%  "
$ws.Range("D30").Value = "'6.90"
$ws.Range("E30").Value = "  +3.68%  "
$ws.Range("D31").Value = "'613.67"
$ws.Range("E31").Value = "  +5.15%  "
$ws.Range("D32").Value = "'3.82"
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("D33").Value = "4.053.10"
$ws.Range("E33").Value = "  +7.24%  "
$ws.Range("D34").Value = "'11.11"
$ws.Range("E34").Value = "  +1.31%  "
$ws.Range("D35").Value = "'0.106"
$ws.Range("E35").Value = "  +2.45%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").Value = "'56.34"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").Value = "'2.79"
$ws.Range("E38").Value = "  +5.32%  "
$ws.Range("E39").Value = "  +4.94%  "
$ws.Range("D40").Value = "'33.96"
$ws.Range("E40").Value = "  -0.95%  "
$ws.Range("D41").Value = "'3.25"
$ws.Range("E41").Value = "  +2.69%  "
$ws.Range("D42").Value = "0.0₃0701"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("E43").Value = "  +0.92%  "
$ws.Range("D44").Value = "'0.343"
$ws.Range("E44").Value = "  +1.54%  "
$ws.Range("D45").Value = "'0.0422"
$ws.Range("E45").Value = "  +2.13%  "
$ws.Range("E46").Value = "  +1.07%  "
$ws.Range("D47").Value = "'2.59"
$ws.Range("E47").Value = "  +1.65%  "
$ws.Range("E48").Value = "  +11.79%  "
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("D50").Value = "'127.50"
$ws.Range("E50").Value = "  +1.26%  "
$ws.Range("D51").Value = "'7.71"
$ws.Range("E51").Value = "  +4.36%  "
